# Standardize spelling and format
#
# On the "Volume By Browser" sheet, columns B/C/D (including the header row)
# were populated in the wrong order: all_sessions, transactions, qty.
# Rotate them one step to the left so the order becomes the intended
# transactions, qty, all_sessions (matching the neighboring "ecr" column,
# which is transactions / qty based).
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Volume By Browser")

$firstRow = $ws.UsedRange.Row
$lastRow = $firstRow + $ws.UsedRange.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldB = $ws.Cells.Item($r, 2).Value2
    $oldC = $ws.Cells.Item($r, 3).Value2
    $oldD = $ws.Cells.Item($r, 4).Value2

    $ws.Cells.Item($r, 2).Value2 = $oldC
    $ws.Cells.Item($r, 3).Value2 = $oldD
    $ws.Cells.Item($r, 4).Value2 = $oldB
}
